$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- New row 19: Datum (date, reuse the date style/format from row 18) ---
$ws.Range("A19").Value = 44185
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats

# --- New row 19: Zeit (hours) ---
$ws.Range("B19").Value = 3

# --- New row 19: Zwischensumme (running total formula) ---
$ws.Range("C19").Formula = "=C18+B19"

# --- New row 19: Tätigkeit (activity text, reuse wrap-text style from row 18) ---
$ws.Range("D19").Value = "Versucht Redpitaya per VPN, über Florians wireguard verfügbar zu machen. --> keine kernel headers für 4.9.0-xilinx und kein module support.."
$ws.Range("D18").Copy()
$ws.Range("D19").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Row height for the new (wrapped, multi-line) row
$ws.Rows.Item(19).RowHeight = 60

# Move the selection to match where Excel leaves it after filling the row
$ws.Range("E19").Select()
